# Updates the cryptos list (Price / Volume(1h) columns, plus a Polkadot /
# ShibaInu row swap) to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold text values (e.g. "409.40", "0.0000219") even
# though they look numeric - force text formatting first so Excel doesn't
# silently coerce them to numbers (and drop trailing zeros / switch to
# scientific notation) when the values are assigned below.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.967.38"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.412.05"
$ws.Range("E3").Value = "  -0.59%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.44%  "

# Row 5 - BNB
$ws.Range("D5").Value = "409.40"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6 - Solana
$ws.Range("D6").Value = "129.42"
$ws.Range("E6").Value = "  -1.11%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +6.34%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.13%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.735"
$ws.Range("E9").Value = "  +6.20%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  +3.08%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "42.86"
$ws.Range("E11").Value = "  +2.10%  "

# Row 12 & 13 - Polkadot and ShibaInu swap places (ShibaInu now ranked 12th,
# Polkadot moves to 13th), with refreshed price/volume figures.
$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D12").Value = "0.0000219"
$ws.Range("E12").Value = "  +43.26%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "9.23"
$ws.Range("E13").Value = "  +9.49%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.22%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.951.80"
$ws.Range("E15").Value = "  -0.53%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "21.25"
$ws.Range("E16").Value = "  +7.32%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.406.46"
$ws.Range("E17").Value = "  -0.84%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "12.55"
$ws.Range("E18").Value = "  +8.26%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  +7.20%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "61.943.13"
$ws.Range("E20").Value = "  -1.24%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "449.58"
$ws.Range("E21").Value = "  +43.72%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "91.78"
$ws.Range("E22").Value = "  +8.78%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.22"
$ws.Range("E23").Value = "  +1.42%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "13.28"
$ws.Range("E24").Value = "  +3.77%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "3.29"
$ws.Range("E25").Value = "  +3.57%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "9.32"
$ws.Range("E26").Value = "  +14.97%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "33.14"
$ws.Range("E27").Value = "  +11.48%  "

# Row 28 - LEO
$ws.Range("D28").Value = "4.80"
$ws.Range("E28").Value = "  +0.70%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  -1.43%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -3.37%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "12.07"
$ws.Range("E31").Value = "  +6.16%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.31%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +0.01%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "42.91"
$ws.Range("E34").Value = "  -3.71%  "

# Row 35 - Dai
$ws.Range("E35").Value = "  -0.19%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0503"
$ws.Range("E36").Value = "  +4.03%  "

# Row 37 - OKB
$ws.Range("D37").Value = "53.90"
$ws.Range("E37").Value = "  +4.10%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.40%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").Value = "3.38"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40 - Stellar
$ws.Range("D40").Value = "0.135"
$ws.Range("E40").Value = "  +7.55%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").Value = "  -0.22%  "

# Row 42 - TheGraph
$ws.Range("D42").Value = "0.319"
$ws.Range("E42").Value = "  -1.14%  "

# Row 43 - Monero
$ws.Range("D43").Value = "142.71"
$ws.Range("E43").Value = "  +0.50%  "

# Row 44 - NEARProtocol
$ws.Range("D44").Value = "4.26"
$ws.Range("E44").Value = "  +8.55%  "

# Row 45 - WEMIXToken
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  +15.85%  "

# Row 46 - ARBITRUM
$ws.Range("E46").Value = "  +1.19%  "

# Row 47 - Celestia
$ws.Range("D47").Value = "16.62"
$ws.Range("E47").Value = "  -1.32%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "0.148"
$ws.Range("E48").Value = "  +23.36%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "22.51"
$ws.Range("E49").Value = "  +6.02%  "

# Row 50 - ThetaToken
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  +8.91%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "3.756.34"
$ws.Range("E51").Value = "  -0.60%  "
